$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder the title/uri rows (columns A and E, rows 4-11) to reflect the
# shuffled order introduced by adding a new JSON record for time bucket
# analysis. Only the displayed text of these cells changes; the underlying
# hyperlinks (by cell reference) are untouched.

$ws.Range("A4").Value = "Ice storm prompts emergency in Oklahoma"
$ws.Range("E4").Value = "https://web.archive.org/web/20080111040510/http://www.cnn.com/2007/WEATHER/01/12/ice.storm/index.html?section=cnn_latest"

$ws.Range("A6").Value = "Snow, wind dump on Colorado"
$ws.Range("E6").Value = "https://web.archive.org/web/20070123183528/http://www.cnn.com/2007/WEATHER/01/21/winter.weather.ap/index.html"

$ws.Range("A7").Value = "Storms kill, knock out power, cancel flights"
$ws.Range("E7").Value = "https://web.archive.org/web/20070122122847/http://www.cnn.com/2007/WEATHER/01/14/winter.weather.ap/index.html"

$ws.Range("A8").Value = "NWS Springfield, MO Homepage"
$ws.Range("E8").Value = "http://www.crh.noaa.gov/sgf/?n=icestormjan07summary%5E"

$ws.Range("A9").Value = "Winter storm blamed for 51 deaths in nine states"
$ws.Range("E9").Value = "https://web.archive.org/web/20070119054443/http://www.cnn.com/2007/WEATHER/01/16/winter.blast.ap/index.html"

$ws.Range("A11").Value = "OzarksFirst.com"
$ws.Range("E11").Value = "https://web.archive.org/web/20070829010112/http://ozarksfirst.com/content/fulltext/?cid=4016"

$wb.Save()
